$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep text formatting so values
# like "312.71" and "1.16%" are stored as text, not numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '312.71'
$ws.Range("E2").Value = '1.16%'

# Row 3
$ws.Range("D3").Value = '37.71'
$ws.Range("E3").Value = '1.17%'

# Row 4
$ws.Range("D4").Value = '5.143'
$ws.Range("E4").Value = '0.55%'

# Row 5
$ws.Range("D5").Value = '0.07913'
$ws.Range("E5").Value = '0.67%'

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.422'
$ws.Range("E6").Value = '1.03%'

# Row 7
$ws.Range("E7").Value = '-0.29%'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '1.901'
$ws.Range("E8").Value = '-3.61%'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '2.976'
$ws.Range("E9").Value = '-4.04%'

# Row 10
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Value = '0.9209'
$ws.Range("E10").Value = '-0.43%'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.1225'
$ws.Range("E11").Value = '-8.38%'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '0.1933'
$ws.Range("E12").Value = '-0.68%'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.09127'
$ws.Range("E13").Value = '2.13%'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03311'
$ws.Range("E14").Value = '-3.90%'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09633'
$ws.Range("E15").Value = '-0.71%'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001378'
$ws.Range("E16").Value = '0.05%'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.005756'
$ws.Range("E17").Value = '-3.35%'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.529'
$ws.Range("E18").Value = '-1.72%'

# Row 19
$ws.Range("E19").Value = '1.55%'

# Row 20
$ws.Range("D20").Value = '5.263'
$ws.Range("E20").Value = '5.13%'

# Row 21
$ws.Range("D21").Value = '0.1273'
$ws.Range("E21").Value = '-1.63%'

# Row 22
$ws.Range("D22").Value = '0.2592'
$ws.Range("E22").Value = '4.13%'

# Row 23
$ws.Range("E23").Value = '-0.19%'

# Row 24
$ws.Range("D24").Value = '0.04363'
$ws.Range("E24").Value = '0.68%'

# Row 25
$ws.Range("D25").Value = '0.001249'
$ws.Range("E25").Value = '2.53%'

# Row 26
$ws.Range("E26").Value = '-4.93%'

# Row 27
$ws.Range("D27").Value = '0.0001220'
$ws.Range("E27").Value = '-9.77%'

# Row 39
$ws.Range("D39").Value = '0.02150'
$ws.Range("E39").Value = '-5.62%'

# Row 40
$ws.Range("D40").Value = '0.05173'
$ws.Range("E40").Value = '3.18%'

# Row 41
$ws.Range("D41").Value = '0.007565'
$ws.Range("E41").Value = '-0.80%'

# Row 42
$ws.Range("E42").Value = '-7.39%'

# Row 43
$ws.Range("D43").Value = '0.1359'
$ws.Range("E43").Value = '0.45%'

# Row 44
$ws.Range("D44").Value = '0.002009'
$ws.Range("E44").Value = '-2.58%'

# Row 45
$ws.Range("D45").Value = '0.008620'
$ws.Range("E45").Value = '2.11%'

# Row 46
$ws.Range("D46").Value = '0.00006715'
$ws.Range("E46").Value = '-0.73%'

# Row 47
$ws.Range("E47").Value = '-0.15%'

# Row 48
$ws.Range("D48").Value = '0.003309'
$ws.Range("E48").Value = '10.12%'

# Row 49
$ws.Range("E49").Value = '-7.74%'

# Row 50
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").Value = '-0.15%'

# Row 51
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").Value = '-0.15%'
